# Edit sheet Card24: append a new, fully-blank row 13 (A13:N13) right
# below the existing data (rows 1-12), extending the used range from
# A1:N12 to A1:N13 - mirrors the upstream diff which adds 14 empty
# inline-string cells on a new row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$lastCol = 14   # column N
$newRow  = 13

$rng = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))

# Writing a lone apostrophe makes Excel store an explicit empty *text*
# value for every cell (rather than leaving them truly untouched, which
# would not extend the worksheet's used range/dimension at all).
$rng.Value = "'"

# Drop the "quote prefix" formatting that the apostrophe entry implies so
# the new cells end up with the same (default) style as their neighbors.
$rng.ClearFormats()

Write-Host "Inserted blank row $newRow (A$newRow`:N$newRow) on sheet 'Card24'"
